$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 1253.3636
$ws.Range("I2").Value = 1473.7142
$ws.Range("J2").Value = 867.75
$ws.Range("K2").Value = 1473.7142
$ws.Range("L2").Value = 867.75
$ws.Range("M2").Value = -1360.7142
$ws.Range("N2").Value = -1093.75

# Row 103
$ws.Range("H103").Value = 333732.66
$ws.Range("I103").Value = 333732.66
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 1001197.98
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -1000611.98
$ws.Range("N103").ClearContents()

# Row 107
$ws.Range("H107").Value = 1214.2174
$ws.Range("I107").Value = 1471
$ws.Range("J107").Value = 732.75
$ws.Range("K107").Value = 1471
$ws.Range("L107").Value = 732.75
$ws.Range("M107").Value = 449
$ws.Range("N107").Value = -4572.75

# Row 116
$ws.Range("H116").Value = 7757.5
$ws.Range("I116").Value = 11265.417
$ws.Range("K116").Value = 11265.417
$ws.Range("M116").Value = -7823.416999999999

# Row 137
$ws.Range("H137").Value = 2253649.2
$ws.Range("I137").Value = 4902902.5
$ws.Range("J137").Value = 1784.2
$ws.Range("K137").Value = 14708707.5
$ws.Range("L137").Value = 5352.6
$ws.Range("M137").Value = -14706157.5
$ws.Range("N137").Value = -10452.6

# Row 140
$ws.Range("H140").Value = 73639.60000000001
$ws.Range("J140").Value = 73639.60000000001
$ws.Range("L140").Value = 73639.60000000001
$ws.Range("N140").Value = -83999.60000000001

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6592964
$ws.Range("I32").Value = 7704641.5
$ws.Range("J32").Value = 23959.182
$ws.Range("K32").Value = 7704641.5
$ws.Range("L32").Value = 23959.182
$ws.Range("M32").Value = -7704354.5
$ws.Range("N32").Value = -24533.182

# Row 44
$ws.Range("H44").Value = 27399.2

# Row 55
$ws.Range("H55").Value = 26249.25

# Row 102
$ws.Range("H102").Value = 3090
$ws.Range("I102").Value = 2612.5
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 2612.5
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = -990.5
$ws.Range("N102").Value = -8244

# Row 108
$ws.Range("H108").Value = 82894.664
$ws.Range("J108").Value = 82894.664
$ws.Range("L108").Value = 82894.664
$ws.Range("N108").Value = -90574.664

# Row 123
$ws.Range("H123").Value = 98429
$ws.Range("J123").Value = 98429
$ws.Range("L123").Value = 98429
$ws.Range("N123").Value = -108229

$ws = $wb.Worksheets.Item("BSM")
# Row 50
$ws.Range("H50").Value = 96890
$ws.Range("J50").Value = 96890
$ws.Range("L50").Value = 96890
$ws.Range("N50").Value = -98038

# Row 86
$ws.Range("H86").Value = 2105
$ws.Range("I86").Value = 2014.25
$ws.Range("J86").Value = 3012.5
$ws.Range("K86").Value = 2014.25
$ws.Range("L86").Value = 3012.5
$ws.Range("M86").Value = -891.25
$ws.Range("N86").Value = -5258.5

# Row 89
$ws.Range("H89").Value = 2105
$ws.Range("I89").Value = 2014.25
$ws.Range("J89").Value = 3012.5
$ws.Range("K89").Value = 10071.25
$ws.Range("L89").Value = 15062.5
$ws.Range("M89").Value = -4455.25
$ws.Range("N89").Value = -26294.5

# Row 96
$ws.Range("H96").Value = 19880
$ws.Range("I96").Value = 11000
$ws.Range("J96").Value = 22100
$ws.Range("K96").Value = 11000
$ws.Range("L96").Value = 22100
$ws.Range("M96").Value = -8254
$ws.Range("N96").Value = -27592

# Row 97
$ws.Range("H97").Value = 20926.143
$ws.Range("I97").Value = 2802.4
$ws.Range("K97").Value = 2802.4
$ws.Range("M97").Value = -1811.4

# Row 110
$ws.Range("H110").Value = 98702
$ws.Range("J110").Value = 98702
$ws.Range("L110").Value = 98702
$ws.Range("N110").Value = -106882

# Row 115
$ws.Range("H115").Value = 76842
$ws.Range("J115").Value = 76842
$ws.Range("L115").Value = 76842
$ws.Range("N115").Value = -79976

# Row 116
$ws.Range("H116").Value = 23296.6
$ws.Range("J116").Value = 23296.6
$ws.Range("L116").Value = 23296.6
$ws.Range("N116").Value = -32474.6

# Row 118
$ws.Range("H118").Value = 45000
$ws.Range("J118").Value = 45000
$ws.Range("L118").Value = 45000
$ws.Range("N118").Value = -48314

# Row 129
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5685.857
$ws.Range("I31").Value = 1689.24
$ws.Range("J31").Value = 7906.2
$ws.Range("K31").Value = 1689.24
$ws.Range("L31").Value = 7906.2
$ws.Range("M31").Value = -1394.24
$ws.Range("N31").Value = -8496.200000000001

# Row 34
$ws.Range("H34").Value = 5685.857
$ws.Range("I34").Value = 1689.24
$ws.Range("J34").Value = 7906.2
$ws.Range("K34").Value = 1689.24
$ws.Range("L34").Value = 7906.2
$ws.Range("M34").Value = -1487.24
$ws.Range("N34").Value = -8310.200000000001

# Row 58
$ws.Range("H58").Value = 1184.0588
$ws.Range("I58").Value = 925.0833
$ws.Range("J58").Value = 1805.6
$ws.Range("K58").Value = 925.0833
$ws.Range("L58").Value = 1805.6
$ws.Range("M58").Value = -722.0833
$ws.Range("N58").Value = -2211.6

# Row 105
$ws.Range("H105").Value = 1385.625
$ws.Range("I105").Value = 1195.5555
$ws.Range("J105").Value = 1630
$ws.Range("K105").Value = 1195.5555
$ws.Range("L105").Value = 1630
$ws.Range("M105").Value = 551.4445000000001
$ws.Range("N105").Value = -5124

# Row 106
$ws.Range("H106").Value = 50000
$ws.Range("J106").Value = 50000
$ws.Range("L106").Value = 50000
$ws.Range("N106").Value = -52524

# Row 117
$ws.Range("H117").Value = 30000
$ws.Range("J117").Value = 30000
$ws.Range("L117").Value = 30000
$ws.Range("N117").Value = -39178

# Row 119
$ws.Range("H119").Value = 38333.332
$ws.Range("J119").Value = 38333.332
$ws.Range("L119").Value = 38333.332
$ws.Range("N119").Value = -48009.332

# Row 127
$ws.Range("H127").Value = 56260
$ws.Range("J127").Value = 56260
$ws.Range("L127").Value = 56260
$ws.Range("N127").Value = -66180

# Row 136
$ws.Range("H136").Value = 1184.0588
$ws.Range("I136").Value = 925.0833
$ws.Range("J136").Value = 1805.6
$ws.Range("K136").Value = 2775.2499
$ws.Range("L136").Value = 5416.799999999999
$ws.Range("M136").Value = -225.2498999999998
$ws.Range("N136").Value = -10516.8

$ws = $wb.Worksheets.Item("CUL")
# Row 58
$ws.Range("H58").Value = 1655.1923

$ws = $wb.Worksheets.Item("GSM")
# Row 62
$ws.Range("H62").Value = 28068.5
$ws.Range("J62").Value = 31076.111
$ws.Range("L62").Value = 31076.111
$ws.Range("N62").Value = -32448.111

# Row 65
$ws.Range("H65").Value = 28068.5
$ws.Range("J65").Value = 31076.111
$ws.Range("L65").Value = 93228.333
$ws.Range("N65").Value = -100092.333

# Row 99
$ws.Range("H99").Value = 9703.556
$ws.Range("J99").Value = 25000
$ws.Range("L99").Value = 25000
$ws.Range("N99").Value = -29492

# Row 108
$ws.Range("H108").Value = 39999.5
$ws.Range("J108").Value = 39999.5
$ws.Range("L108").Value = 39999.5
$ws.Range("N108").Value = -47679.5

$ws = $wb.Worksheets.Item("LTW")
# Row 70
$ws.Range("H70").Value = 35400
$ws.Range("J70").Value = 35400
$ws.Range("L70").Value = 35400
$ws.Range("N70").Value = -35940

# Row 73
$ws.Range("H73").Value = 35400
$ws.Range("J73").Value = 35400
$ws.Range("L73").Value = 35400
$ws.Range("N73").Value = -37272

# Row 93
$ws.Range("H93").Value = 51500
$ws.Range("I93").Value = 100000
$ws.Range("J93").Value = 3000
$ws.Range("K93").Value = 100000
$ws.Range("L93").Value = 3000
$ws.Range("M93").Value = -98752
$ws.Range("N93").Value = -5496

# Row 123
$ws.Range("H123").Value = 40000
$ws.Range("J123").Value = 40000
$ws.Range("L123").Value = 40000
$ws.Range("N123").Value = -49800

$ws = $wb.Worksheets.Item("WVR")
# Row 95
$ws.Range("H95").Value = 95000
$ws.Range("J95").Value = 95000
$ws.Range("L95").Value = 95000
$ws.Range("N95").Value = -100492
